$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 119; this shifts existing rows 119-171 down to 120-172
# (Excel automatically carries formatting/styles along with the shift).
$ws.Rows(119).EntireRow.Insert()

# Populate the newly inserted row 119 with the new record's data.
$ws.Range("A119").Value = 1
$ws.Range("B119").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C119").Value = "Arica y Parinacota"
$ws.Range("D119").Value = 45009
$ws.Range("E119").Value = 15
$ws.Range("F119").Value = "Fruta"
$ws.Range("G119").Value = 100106
$ws.Range("H119").Value = "Oleaginosos"
$ws.Range("I119").Value = 100106002
$ws.Range("J119").Value = "Palta"
$ws.Range("K119").Value = "Hass"
$ws.Range("L119").Value = "Segunda"
$ws.Range("M119").Value = 250
$ws.Range("N119").Value = 25000
$ws.Range("O119").Value = 26000
$ws.Range("P119").Value = 25600
$ws.Range("Q119").Value = "$/bandeja 10 kilos"
$ws.Range("R119").Value = "Perú"
$ws.Range("S119").Value = 2560
$ws.Range("T119").Value = 10
